# Update table from Cruise to Month:
#  - remove the "Habitat" column (old column B)
#  - rename "Cruise" header to "Month"
#  - replace cruise codes with the month the cruise took place

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map old cruise code (still sitting in column A before the shift) to month name.
$cruiseToMonth = @{
    "OR1-1219" = "March"
    "OR1-1242" = "October"
}

# Replace column A values (cruise code -> month) before the column shift so we
# work off of the pre-edit layout.
for ($r = 2; $r -le 12; $r++) {
    $cruise = $ws.Cells.Item($r, 1).Value()
    if ($cruiseToMonth.ContainsKey($cruise)) {
        $ws.Cells.Item($r, 1).Value = $cruiseToMonth[$cruise]
    }
}

# Rename the header in column A from "Cruise" to "Month".
$ws.Range("A1").Value = "Month"

# Delete the whole "Habitat" column (column B) - this shifts every column
# after it one to the left, matching the diff (Station moves from C to B,
# Date from D to C, etc.) automatically.
$ws.Columns.Item(2).Delete()
